# This script reproduces the commit:
#   "Fruta / hortaliza, semanal"
# which adds one new weekly price-record row to the dataset.
#
# The new record is inserted as row 235 (pushing the existing rows
# 235-306 down to 236-307), so the sheet dimension grows from
# A1:R306 to A1:R307.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 235; Excel automatically shifts
# every row at/after 235 down by one (preserving their values/styles),
# exactly matching the target workbook layout.
$ws.Rows("235:235").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(235, 1).Value  = 6
$ws.Cells.Item(235, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(235, 3).Value  = "Metropolitana"
$ws.Cells.Item(235, 4).Value  = 44543
$ws.Cells.Item(235, 5).Value  = 13
$ws.Cells.Item(235, 6).Value  = 100112030
$ws.Cells.Item(235, 7).Value  = "Poroto granado"
$ws.Cells.Item(235, 8).Value  = "Sin especificar"
$ws.Cells.Item(235, 9).Value  = "Primera"
$ws.Cells.Item(235, 10).Value = 200
$ws.Cells.Item(235, 11).Value = 30000
$ws.Cells.Item(235, 12).Value = 35000
$ws.Cells.Item(235, 13).Value = 33000
$ws.Cells.Item(235, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(235, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(235, 16).Value = 2200
$ws.Cells.Item(235, 17).Value = 15
$ws.Cells.Item(235, 18).Value = "Hortaliza"
